$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "(according to the population census data)" note row (old row 2).
# This shifts the "(sq. km)" row, the years header row, and the Area data
# row up by one, and drops the now-unused shared string.
$ws.Rows("2:2").Delete()

# Only the most recent year (2014 / 1082.5, originally in column D) should
# remain; drop the 1989 and 2002 columns (B and C) so column D becomes B.
$ws.Columns("B:C").Delete()

# Re-apply a uniform row height to the whole table (including a trailing
# blank row) to match the re-exported layout.
$ws.Range("A1:B6").RowHeight = 20.1
